$wb = $excel.ActiveWorkbook

# --- STREAMS sheet updates ---
$ws = $wb.Worksheets.Item("STREAMS")

# Update existing cell: liquid flow rate (real/nonideal) for the first stream
$ws.Range("L4").Value = 518.1

# Add a new row (27) for the "Liquid flow rate nonideal" (Lreal) variable.
# Copy formatting from the row above (row 26) first so the new row matches
# the existing table styling, then set the cell contents.
$ws.Range("A26:N26").Copy($ws.Range("A27:N27"))

$ws.Range("A27").Value = "Liquid flow rate nonideal"
$ws.Range("B27").Value = "Lreal"
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 575.6
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = "mol.s-1"

# Select the newly added row's L cell and make STREAMS the active sheet,
# matching the saved workbook view state.
[void]$ws.Range("L27").Select()
$ws.Activate()
